# Apply cryptos list update (Fri May 17 02:18:11 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.214.27"
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.942.78"
$ws.Range("E3").Value = "  -2.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.34"
$ws.Range("E5").Value = "  -2.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.78"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  -0.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.939.81"
$ws.Range("E9").Value = "  -2.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  -4.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  -3.42%  "

$ws.Range("E12").Value = "  +0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.43"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.266.96"
$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.408.76"
$ws.Range("E17").Value = "  -3.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.05"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.942.81"
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.76"
$ws.Range("E20").Value = "  +13.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.36"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  -0.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.52"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.16"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -6.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.45"
$ws.Range("E29").Value = "  +3.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.98"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0999"
$ws.Range("E32").Value = "  -5.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.17"
$ws.Range("E34").Value = "  +0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.71"
$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("B38").Value = "Arweave"
$ws.Range("C38").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "45.14"
$ws.Range("E38").Value = "  +2.99%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.19"
$ws.Range("E39").Value = "  -1.82%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  -8.87%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.301"
$ws.Range("E41").Value = "  -1.46%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.54"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  -6.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "382.27"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.677.04"
$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.58"
$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.17"
$ws.Range("E50").Value = "  +1.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.53"
$ws.Range("E51").Value = "  -1.60%  "
